$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns I0/IF - copy format from existing header (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows (row, I value, J value)
$data = @(
    "2,4,5"
    "3,7,7"
    "4,6,7"
    "5,6,7"
    "6,7,8"
    "7,6,7"
    "8,7,7"
    "9,7,8"
    "10,8,8"
    "11,6,6"
    "12,8,8"
    "13,9,9"
    "14,6,6"
    "15,8,8"
    "16,7,7"
    "17,7,7"
    "18,6,6"
    "19,8,8"
    "20,7,7"
    "21,7,8"
    "22,8,8"
    "23,7,7"
    "24,7,7"
    "25,7,7"
    "26,8,8"
    "27,8,8"
    "28,9,9"
    "29,7,7"
    "30,8,8"
    "31,7,7"
    "32,7,7"
    "33,7,7"
    "34,8,8"
    "35,8,8"
    "36,8,8"
    "37,7,7"
    "38,8,8"
    "39,7,7"
    "40,8,8"
    "41,8,8"
    "42,7,7"
    "43,7,7"
    "44,7,7"
    "45,7,7"
    "46,8,8"
    "47,7,7"
    "48,8,8"
    "49,7,7"
    "50,8,8"
    "51,8,8"
    "52,8,8"
    "53,8,8"
    "54,7,7"
    "55,7,7"
    "56,6,6"
    "57,8,8"
    "58,8,8"
    "59,8,8"
    "60,8,8"
    "61,8,8"
    "62,9,9"
    "63,7,7"
    "64,7,7"
    "65,7,7"
    "66,7,7"
    "67,7,7"
    "68,8,8"
    "69,6,6"
    "70,8,8"
    "71,8,8"
    "72,8,8"
    "73,7,7"
    "74,7,7"
    "75,7,7"
    "76,8,8"
    "77,10,11"
    "78,8,9"
    "79,7,7"
    "80,5,5"
    "81,6,6"
    "82,4,4"
    "83,5,5"
    "84,5,5"
    "85,4,4"
)

foreach ($entry in $data) {
    $parts = $entry.Split(",")
    $r = [int]$parts[0]
    $iVal = [int]$parts[1]
    $jVal = [int]$parts[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
